$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value, taken from the latest coinranking.com scrape.
# Cells whose new text is a "pure" number (e.g. "1.000", "0.9989") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# and silently drops the formatting (trailing zeros, leading zeros, etc.).
$updates = [ordered]@{
    'D2' = '29.350.38'
    'E2' = '  -0.50%  '
    'D3' = '1.844.57'
    'E3' = '  -0.27%  '
    'D4' = '0.9989'
    'E4' = '  +0.05%  '
    'D5' = '240.18'
    'E5' = '  -0.67%  '
    'D6' = '0.6311'
    'E6' = '  +0.18%  '
    'D8' = '0.07533'
    'E8' = '  -0.13%  '
    'E9' = '  -0.77%  '
    'D10' = '24.41'
    'E10' = '  +0.00%  '
    'D11' = '0.07710'
    'E11' = '  +0.16%  '
    'D12' = '1.853.45'
    'E12' = '  -1.76%  '
    'D13' = '4.988'
    'E13' = '  -0.47%  '
    'D14' = '0.6822'
    'E14' = '  -0.82%  '
    'D15' = '0.00001001'
    'E15' = '  +2.12%  '
    'D16' = '82.75'
    'E16' = '  -1.23%  '
    'B17' = 'WrappedliquidstakedEther2.0'
    'C17' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D17' = '2.105.92'
    'E17' = '  -3.08%  '
    'B18' = 'Uniswap'
    'C18' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D18' = '6.119'
    'E18' = '  -2.18%  '
    'B19' = 'WrappedBTC'
    'C19' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D19' = '29.387.27'
    'E19' = '  -0.59%  '
    'B20' = 'BitcoinCash'
    'C20' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D20' = '227.59'
    'E20' = '  -2.78%  '
    'B21' = 'Avalanche'
    'C21' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D21' = '12.42'
    'E21' = '  -0.71%  '
    'B22' = 'Dai'
    'C22' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'D22' = '1.000'
    'E22' = '  +0.01%  '
    'B23' = 'Chainlink'
    'C23' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'D23' = '7.533'
    'E23' = '  -1.33%  '
    'B24' = 'BinanceUSD'
    'C24' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D24' = '1.000'
    'E24' = '  +0.10%  '
    'D25' = '157.30'
    'E25' = '  +1.51%  '
    'E26' = '  +0.13%  '
    'D27' = '8.335'
    'E27' = '  -1.38%  '
    'E28' = '  -0.52%  '
    'D29' = '1.466'
    'E29' = '  -0.88%  '
    'E30' = '  -0.29%  '
    'D31' = '0.05669'
    'E31' = '  -3.16%  '
    'D32' = '4.123'
    'E32' = '  +0.46%  '
    'E33' = '  -0.73%  '
    'E34' = '  -2.65%  '
    'E35' = '  -1.49%  '
    'D36' = '0.7128'
    'E36' = '  -1.49%  '
    'E37' = '  +0.17%  '
    'D38' = '1.259.18'
    'E38' = '  +1.25%  '
    'D39' = '0.01807'
    'E39' = '  +1.09%  '
    'D40' = '2.780'
    'E40' = '  -0.55%  '
    'B41' = 'FraxShare'
    'C41' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D41' = '6.221'
    'E41' = '  +1.06%  '
    'B42' = 'TrustWalletToken'
    'C42' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D42' = '0.9090'
    'E42' = '  -0.03%  '
    'D43' = '1.0000'
    'E43' = '  +0.06%  '
    'D44' = '101.16'
    'E44' = '  -0.83%  '
    'D45' = '66.13'
    'E45' = '  -1.71%  '
    'B46' = 'TheSandbox'
    'C46' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D46' = '0.4036'
    'E46' = '  -0.19%  '
    'B47' = 'Aptos'
    'C47' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D47' = '7.039'
    'E47' = '  -3.93%  '
    'B48' = 'BabyDogeCoin'
    'C48' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D48' = '0.00000000117'
    'E48' = '  +0.80%  '
    'D49' = '9.114'
    'E49' = '  -0.75%  '
    'D50' = '1.683'
    'E50' = '  -1.96%  '
    'D51' = '0.1123'
    'E51' = '  +0.41%  '
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $cell = $ws.Range($ref)
    if ($value -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        # Numeric-looking text (e.g. "1.000", "0.9989") - keep as text
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}
